# Update the Mann-Whitney test results: new alpha_MW p-values and
# recomputed "significant" flags (col C) for the fish prey/not-prey nutrient comparison.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: As
$ws.Range("B2").Value = 0.000712758942006847
$ws.Range("C2").Value = "yes"

# Row 3: Ca
$ws.Range("B3").Value = 0.00000245083453239043
$ws.Range("C3").Value = "yes"

# Row 4: Co
$ws.Range("B4").Value = 0.000000437603813112014
$ws.Range("C4").Value = "yes"

# Row 5: Cu
$ws.Range("B5").Value = 0.000188909683955304
$ws.Range("C5").Value = "yes"

# Row 6: Fe
$ws.Range("B6").Value = 0.0000773285943019888
$ws.Range("C6").Value = "yes"

# Row 7: K
$ws.Range("B7").Value = 0.0186279809634519
$ws.Range("C7").Value = "yes"

# Row 8: Mg
$ws.Range("B8").Value = 0.0816556651137754
$ws.Range("C8").Value = "no"

# Row 9: Mn
$ws.Range("B9").Value = 0.0000847274067207474
$ws.Range("C9").Value = "yes"

# Row 10: Na
$ws.Range("B10").Value = 0.000000000171555091625599
$ws.Range("C10").Value = "yes"

# Row 11: Ni
$ws.Range("B11").Value = 0.000000159070685564544
$ws.Range("C11").Value = "yes"

# Row 12: P
$ws.Range("B12").Value = 0.000000000125656593430103
$ws.Range("C12").Value = "yes"

# Row 13: Se
$ws.Range("B13").Value = 0.000000000000194275607585403
$ws.Range("C13").Value = "yes"

# Row 14: Zn
$ws.Range("B14").Value = 0.0000546708964788386
$ws.Range("C14").Value = "yes"
